# Add data for 2025-12-25
# Applies the verified cell-level deltas for the new day across the
# Citywide Totals, By Neighborhood rollup, and individual neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 6490
$ws.Range("L3").Value = 6997
$ws.Range("K4").Value = 1790
$ws.Range("L4").Value = 1744
$ws.Range("L5").Value = 413
$ws.Range("L6").Value = 5748
$ws.Range("K7").Value = 27584
$ws.Range("L7").Value = 21392

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L2").Value = 191
$ws.Range("L7").Value = 681
$ws.Range("L8").Value = 1412
$ws.Range("L11").Value = 351
$ws.Range("L14").Value = 103
$ws.Range("L15").Value = 184
$ws.Range("L18").Value = 147
$ws.Range("L19").Value = 589
$ws.Range("K20").Value = 668
$ws.Range("L20").Value = 539
$ws.Range("L24").Value = 66
$ws.Range("L29").Value = 1195
$ws.Range("L33").Value = 962
$ws.Range("L34").Value = 117
$ws.Range("L37").Value = 823
$ws.Range("L42").Value = 676
$ws.Range("L43").Value = 160
$ws.Range("L44").Value = 147
$ws.Range("L50").Value = 104
$ws.Range("L51").Value = 264
$ws.Range("L52").Value = 455
$ws.Range("L63").Value = 65
$ws.Range("L64").Value = 136
$ws.Range("L65").Value = 424
$ws.Range("L67").Value = 744
$ws.Range("L71").Value = 54
$ws.Range("L75").Value = 78
$ws.Range("L76").Value = 339
$ws.Range("L79").Value = 593
$ws.Range("L83").Value = 470
$ws.Range("L85").Value = 1063
$ws.Range("L91").Value = 288
$ws.Range("L94").Value = 260
$ws.Range("L95").Value = 297
$ws.Range("L96").Value = 235
$ws.Range("L99").Value = 368
$ws.Range("K101").Value = 27584
$ws.Range("L101").Value = 21392

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("L2").Value = 44
$ws.Range("L7").Value = 103

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L2").Value = 76
$ws.Range("L7").Value = 235

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 237
$ws.Range("L3").Value = 219
$ws.Range("L7").Value = 681

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L2").Value = 133
$ws.Range("L6").Value = 90
$ws.Range("L7").Value = 351

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 317
$ws.Range("L3").Value = 440
$ws.Range("L4").Value = 62
$ws.Range("L6").Value = 223
$ws.Range("L7").Value = 1063

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L2").Value = 145
$ws.Range("L3").Value = 142
$ws.Range("L7").Value = 455

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L4").Value = 99
$ws.Range("L6").Value = 339
$ws.Range("L7").Value = 1412

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L4").Value = 19
$ws.Range("L6").Value = 106
$ws.Range("L7").Value = 470

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L3").Value = 340
$ws.Range("L4").Value = 66
$ws.Range("L6").Value = 270
$ws.Range("L7").Value = 962

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L2").Value = 106
$ws.Range("L7").Value = 297

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 250
$ws.Range("L3").Value = 292
$ws.Range("L7").Value = 823

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L2").Value = 158
$ws.Range("L6").Value = 104
$ws.Range("L7").Value = 424

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L3").Value = 147
$ws.Range("L7").Value = 368

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 211
$ws.Range("L7").Value = 744

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 365
$ws.Range("L5").Value = 20
$ws.Range("L7").Value = 1195

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 211
$ws.Range("L7").Value = 589

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L2").Value = 57
$ws.Range("L3").Value = 42
$ws.Range("L7").Value = 147

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L3").Value = 67
$ws.Range("L6").Value = 153
$ws.Range("L7").Value = 339

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L3").Value = 233
$ws.Range("L6").Value = 192
$ws.Range("L7").Value = 676

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("L2").Value = 26
$ws.Range("L7").Value = 66

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L3").Value = 129
$ws.Range("L7").Value = 288

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L2").Value = 184
$ws.Range("L6").Value = 158
$ws.Range("L7").Value = 593

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("L3").Value = 42
$ws.Range("L7").Value = 136

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L3").Value = 188
$ws.Range("K4").Value = 32
$ws.Range("K7").Value = 668
$ws.Range("L7").Value = 539

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("L3").Value = 51
$ws.Range("L7").Value = 147

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("L3").Value = 32
$ws.Range("L7").Value = 117

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L6").Value = 96
$ws.Range("L7").Value = 260

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L2").Value = 68
$ws.Range("L3").Value = 58
$ws.Range("L7").Value = 184

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("L2").Value = 34
$ws.Range("L7").Value = 104

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("L3").Value = 63
$ws.Range("L4").Value = 16
$ws.Range("L7").Value = 191

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("L3").Value = 28
$ws.Range("L7").Value = 78

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L2").Value = 83
$ws.Range("L7").Value = 264

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("L6").Value = 49
$ws.Range("L7").Value = 160

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("L3").Value = 18
$ws.Range("L7").Value = 54

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("L4").Value = 7
$ws.Range("L6").Value = 34
